$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-19 Sunday" "2025-01-20 Monday"

Replace-Text "183÷3=" "246÷9="
Replace-Text "311÷7=" "264÷7="
Replace-Text "595÷7=" "417÷4="
Replace-Text "446÷5=" "735÷8="
Replace-Text "793÷6=" "501÷7="
Replace-Text "592÷5=" "793÷5="
Replace-Text "315÷4=" "813÷4="
Replace-Text "582÷7=" "892÷7="
Replace-Text "849÷2=" "303÷4="
Replace-Text "880÷2=" "377÷8="
Replace-Text "703÷5=" "494÷8="
Replace-Text "753÷8=" "526÷7="
Replace-Text "115÷4=" "673÷2="
Replace-Text "614÷3=" "138÷4="
Replace-Text "821÷7=" "554÷2="
Replace-Text "865÷6=" "807÷3="
Replace-Text "613÷6=" "797÷7="
Replace-Text "254÷3=" "166÷5="
Replace-Text "890÷7=" "857÷7="
Replace-Text "732÷3=" "379÷5="
Replace-Text "731÷7=" "333÷5="
Replace-Text "515÷2=" "721÷4="
Replace-Text "606÷5=" "370÷3="
Replace-Text "714÷8=" "753÷7="
Replace-Text "463÷7=" "108÷2="

Write-Output "Done"
